# Auto-generated Excel COM-interop script
# Applies numeric cell-value corrections across multiple sheets,
# as produced by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 714373.1
$ws.Cells.Item(17, 10).Value = 714373.1
$ws.Cells.Item(17, 12).Value = 2143119.3
$ws.Cells.Item(17, 14).Value = -2143455.3
$ws.Cells.Item(33, 8).Value = 3571839.8
$ws.Cells.Item(33, 9).Value = 3846555.8
$ws.Cells.Item(33, 11).Value = 3846555.8
$ws.Cells.Item(33, 13).Value = -3846326.8
$ws.Cells.Item(76, 8).Value = 8399.200000000001
$ws.Cells.Item(76, 9).Value = 5999.5
$ws.Cells.Item(76, 10).Value = 9999
$ws.Cells.Item(76, 11).Value = 5999.5
$ws.Cells.Item(76, 12).Value = 9999
$ws.Cells.Item(76, 13).Value = -5684.5
$ws.Cells.Item(76, 14).Value = -10629
$ws.Cells.Item(79, 8).Value = 8399.200000000001
$ws.Cells.Item(79, 9).Value = 5999.5
$ws.Cells.Item(79, 10).Value = 9999
$ws.Cells.Item(79, 11).Value = 5999.5
$ws.Cells.Item(79, 12).Value = 9999
$ws.Cells.Item(79, 13).Value = -4907.5
$ws.Cells.Item(79, 14).Value = -12183
$ws.Cells.Item(106, 8).Value = 4329.6665
$ws.Cells.Item(106, 9).Value = 4329.6665
$ws.Cells.Item(106, 11).Value = 4329.6665
$ws.Cells.Item(106, 13).Value = -3698.6665
$ws.Cells.Item(137, 8).Value = 1490.8572
$ws.Cells.Item(137, 10).Value = 1746.5
$ws.Cells.Item(137, 12).Value = 5239.5
$ws.Cells.Item(137, 14).Value = -10339.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4035.36
$ws.Cells.Item(32, 9).Value = 3527.6956
$ws.Cells.Item(32, 11).Value = 3527.6956
$ws.Cells.Item(32, 13).Value = -3240.6956
$ws.Cells.Item(53, 8).Value = 18255.143
$ws.Cells.Item(53, 9).Value = 16557.2
$ws.Cells.Item(53, 11).Value = 16557.2
$ws.Cells.Item(53, 13).Value = -15875.2
$ws.Cells.Item(61, 8).Value = 2851.9773
$ws.Cells.Item(61, 9).Value = 2430.0303
$ws.Cells.Item(61, 10).Value = 4117.8184
$ws.Cells.Item(61, 11).Value = 2430.0303
$ws.Cells.Item(61, 12).Value = 4117.8184
$ws.Cells.Item(61, 13).Value = -2218.0303
$ws.Cells.Item(61, 14).Value = -4541.8184
$ws.Cells.Item(97, 8).Value = 483.24
$ws.Cells.Item(97, 9).Value = 451.34784
$ws.Cells.Item(97, 11).Value = 451.34784
$ws.Cells.Item(97, 13).Value = 44.65215999999998
$ws.Cells.Item(132, 8).Value = 3039.7896
$ws.Cells.Item(132, 9).Value = 2838.4614
$ws.Cells.Item(132, 10).Value = 3476
$ws.Cells.Item(132, 11).Value = 8515.3842
$ws.Cells.Item(132, 12).Value = 10428
$ws.Cells.Item(132, 13).Value = -5985.3842
$ws.Cells.Item(132, 14).Value = -15488
$ws.Cells.Item(136, 8).Value = 2851.9773
$ws.Cells.Item(136, 9).Value = 2430.0303
$ws.Cells.Item(136, 10).Value = 4117.8184
$ws.Cells.Item(136, 11).Value = 7290.090899999999
$ws.Cells.Item(136, 12).Value = 12353.4552
$ws.Cells.Item(136, 13).Value = -4740.090899999999
$ws.Cells.Item(136, 14).Value = -17453.4552

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(29, 8).Value = 14314.25
$ws.Cells.Item(29, 9).Value = 14930.571
$ws.Cells.Item(29, 11).Value = 14930.571
$ws.Cells.Item(29, 13).Value = -14641.571
$ws.Cells.Item(134, 8).Value = 2268.7273
$ws.Cells.Item(134, 9).Value = 2277.125
$ws.Cells.Item(134, 11).Value = 6831.375
$ws.Cells.Item(134, 13).Value = -4296.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 26209.453
$ws.Cells.Item(31, 9).Value = 35159.3
$ws.Cells.Item(31, 11).Value = 35159.3
$ws.Cells.Item(31, 13).Value = -34864.3
$ws.Cells.Item(34, 8).Value = 26209.453
$ws.Cells.Item(34, 9).Value = 35159.3
$ws.Cells.Item(34, 11).Value = 35159.3
$ws.Cells.Item(34, 13).Value = -34957.3
$ws.Cells.Item(122, 8).Value = 1911.25
$ws.Cells.Item(122, 9).Value = 1826.2222
$ws.Cells.Item(122, 11).Value = 5478.6666
$ws.Cells.Item(122, 13).Value = -3028.6666
$ws.Cells.Item(132, 8).Value = 8647.143
$ws.Cells.Item(132, 9).Value = 5172.1333
$ws.Cells.Item(132, 10).Value = 17334.666
$ws.Cells.Item(132, 11).Value = 15516.3999
$ws.Cells.Item(132, 12).Value = 52003.99800000001
$ws.Cells.Item(132, 13).Value = -12986.3999
$ws.Cells.Item(132, 14).Value = -57063.99800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 39090772
$ws.Cells.Item(4, 9).Value = 41930828
$ws.Cells.Item(4, 11).Value = 125792484
$ws.Cells.Item(4, 13).Value = -125792372
$ws.Cells.Item(20, 8).Value = 1448.3334
$ws.Cells.Item(20, 9).Value = 1345
$ws.Cells.Item(20, 10).Value = 1500
$ws.Cells.Item(20, 11).Value = 4035
$ws.Cells.Item(20, 12).Value = 4500
$ws.Cells.Item(20, 13).Value = -3808
$ws.Cells.Item(20, 14).Value = -4954
$ws.Cells.Item(29, 8).Value = 2143.261
$ws.Cells.Item(29, 9).Value = 306.4
$ws.Cells.Item(29, 10).Value = 5587.375
$ws.Cells.Item(29, 11).Value = 919.1999999999999
$ws.Cells.Item(29, 12).Value = 16762.125
$ws.Cells.Item(29, 13).Value = -642.1999999999999
$ws.Cells.Item(29, 14).Value = -17316.125
$ws.Cells.Item(94, 8).Value = 6131.625
$ws.Cells.Item(94, 9).Value = 3561.2856
$ws.Cells.Item(94, 10).Value = 8130.778
$ws.Cells.Item(94, 11).Value = 10683.8568
$ws.Cells.Item(94, 12).Value = 24392.334
$ws.Cells.Item(94, 13).Value = -10007.8568
$ws.Cells.Item(94, 14).Value = -25744.334
$ws.Cells.Item(108, 8).Value = 793.4
$ws.Cells.Item(108, 9).Value = 793.4
$ws.Cells.Item(108, 11).Value = 2380.2
$ws.Cells.Item(108, 13).Value = 499.8000000000002
$ws.Cells.Item(113, 8).Value = 731.4375
$ws.Cells.Item(113, 9).Value = 432.5
$ws.Cells.Item(113, 11).Value = 1297.5
$ws.Cells.Item(113, 13).Value = 872.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2303.7144
$ws.Cells.Item(80, 9).Value = 2304.3333
$ws.Cells.Item(80, 10).Value = 2300
$ws.Cells.Item(80, 11).Value = 2304.3333
$ws.Cells.Item(80, 12).Value = 2300
$ws.Cells.Item(80, 13).Value = -1306.3333
$ws.Cells.Item(80, 14).Value = -4296
$ws.Cells.Item(83, 8).Value = 2303.7144
$ws.Cells.Item(83, 9).Value = 2304.3333
$ws.Cells.Item(83, 10).Value = 2300
$ws.Cells.Item(83, 11).Value = 11521.6665
$ws.Cells.Item(83, 12).Value = 11500
$ws.Cells.Item(83, 13).Value = -6529.666499999999
$ws.Cells.Item(83, 14).Value = -21484

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4892.6924
$ws.Cells.Item(132, 9).Value = 5066
$ws.Cells.Item(132, 11).Value = 15198
$ws.Cells.Item(132, 13).Value = -12668
$ws.Cells.Item(136, 8).Value = 6600
$ws.Cells.Item(136, 10).Value = 7500
$ws.Cells.Item(136, 12).Value = 22500
$ws.Cells.Item(136, 14).Value = -27600

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2837.7942
$ws.Cells.Item(122, 9).Value = 2754.6897
$ws.Cells.Item(122, 10).Value = 3319.8
$ws.Cells.Item(122, 11).Value = 8264.069100000001
$ws.Cells.Item(122, 12).Value = 9959.400000000001
$ws.Cells.Item(122, 13).Value = -5814.069100000001
$ws.Cells.Item(122, 14).Value = -14859.4
$ws.Cells.Item(132, 8).Value = 2466
$ws.Cells.Item(132, 9).Value = 2449
$ws.Cells.Item(132, 11).Value = 7347
$ws.Cells.Item(132, 13).Value = -4817
$ws.Cells.Item(136, 8).Value = 1967.2667
$ws.Cells.Item(136, 9).Value = 1684.6
$ws.Cells.Item(136, 10).Value = 2532.6
$ws.Cells.Item(136, 11).Value = 5053.799999999999
$ws.Cells.Item(136, 12).Value = 7597.799999999999
$ws.Cells.Item(136, 13).Value = -2503.799999999999
